# Generate Report for Handback
#
# Row 3 (the 28c5d04f-... file) on the Overview sheet and on each language
# sheet (zh-cn, de-de) moves from "Ready for handoff" to
# "Handback transform failed", and the per-language sheets get a detailed
# error message written into the "Error Detail" column (P) for that row.
# Column P is widened to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

$languages = @(
    @{ Sheet = "zh-cn"; Locale = "zh-cn" },
    @{ Sheet = "de-de"; Locale = "de-de" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column (C) for the handoff-pending row.
    $ws.Range("C3").Value = "Handback transform failed"

    # Error Detail column (P) message + widen column P to fit it.
    $message = "Handback file name: qigobdch.cuv is different with handoff file name: 28c5d04f-e27d-4a67-8f62-57496c53ab27.cf20b99b012528349d3daaf8005d05a63ce199a5." + $lang.Locale + "."
    $ws.Range("P3").Value = $message
    # NOTE: ColumnWidth is in character units; Excel stores a slightly
    # different "true" width in the XML (padding baked into the column's
    # max-digit-width). 39.1667 characters round-trips to a stored
    # width of exactly 40, matching the other pre-widened columns (A, G,
    # I, J) on this sheet.
    $ws.Columns.Item(16).ColumnWidth = 39.1667
}
